$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns (D, E, F) before the existing "Terms Typically Offered" column,
# which shifts it from D to G.
$ws.Range("D1:F1").EntireColumn.Insert()

# New header labels for the inserted columns plus the (now shifted) last header.
$ws.Cells.Item(1, 4).Value = "Corequisites"
$ws.Cells.Item(1, 5).Value = "Concurrent"
$ws.Cells.Item(1, 6).Value = "Recommended"

# Per-row data: split any "Corequisite:"/"Concurrent:" clause that used to be appended
# to column C (Prerequisites) out into the new Corequisites / Concurrent columns, and
# backfill "NA" for rows/columns that do not apply. Column G keeps the original
# "Terms Typically Offered" value (already shifted there by the column insert above).

$ws.Cells.Item(2, 4).Value = "NA"
$ws.Cells.Item(2, 5).Value = "NA"
$ws.Cells.Item(2, 6).Value = "NA"
$ws.Cells.Item(2, 7).Value = "F, W, SP"
$ws.Cells.Item(3, 4).Value = "NA"
$ws.Cells.Item(3, 5).Value = "NA"
$ws.Cells.Item(3, 6).Value = "NA"
$ws.Cells.Item(3, 7).Value = "TBD"
$ws.Cells.Item(4, 4).Value = "NA"
$ws.Cells.Item(4, 5).Value = "NA"
$ws.Cells.Item(4, 6).Value = "NA"
$ws.Cells.Item(4, 7).Value = "F"
$ws.Cells.Item(5, 4).Value = "NA"
$ws.Cells.Item(5, 5).Value = "NA"
$ws.Cells.Item(5, 6).Value = "NA"
$ws.Cells.Item(5, 7).Value = "W"
$ws.Cells.Item(6, 4).Value = "NA"
$ws.Cells.Item(6, 5).Value = "NA"
$ws.Cells.Item(6, 6).Value = "NA"
$ws.Cells.Item(6, 7).Value = "SP"
$ws.Cells.Item(7, 3).Value = "ARCH 242."
$ws.Cells.Item(7, 4).Value = "NA"
$ws.Cells.Item(7, 5).Value = "ARCH 253."
$ws.Cells.Item(7, 6).Value = "NA"
$ws.Cells.Item(7, 7).Value = "SP "
$ws.Cells.Item(8, 4).Value = "NA"
$ws.Cells.Item(8, 5).Value = "NA"
$ws.Cells.Item(8, 6).Value = "NA"
$ws.Cells.Item(8, 7).Value = "F"
$ws.Cells.Item(9, 4).Value = "NA"
$ws.Cells.Item(9, 5).Value = "NA"
$ws.Cells.Item(9, 6).Value = "NA"
$ws.Cells.Item(9, 7).Value = "W"
$ws.Cells.Item(10, 4).Value = "NA"
$ws.Cells.Item(10, 5).Value = "NA"
$ws.Cells.Item(10, 6).Value = "NA"
$ws.Cells.Item(10, 7).Value = "SP, SU"
$ws.Cells.Item(11, 3).Value = "ARCH 133."
$ws.Cells.Item(11, 4).Value = "ARCH 251."
$ws.Cells.Item(11, 5).Value = "NA"
$ws.Cells.Item(11, 6).Value = "NA"
$ws.Cells.Item(11, 7).Value = "F "
$ws.Cells.Item(12, 3).Value = "ARCH 241."
$ws.Cells.Item(12, 4).Value = "ARCH 252."
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(12, 6).Value = "NA"
$ws.Cells.Item(12, 7).Value = "W "
$ws.Cells.Item(13, 3).Value = "ARCH 133."
$ws.Cells.Item(13, 4).Value = "ARCH 241."
$ws.Cells.Item(13, 5).Value = "NA"
$ws.Cells.Item(13, 6).Value = "NA"
$ws.Cells.Item(13, 7).Value = "F "
$ws.Cells.Item(14, 3).Value = "ARCH 251, ARCH 241."
$ws.Cells.Item(14, 4).Value = "ARCH 242."
$ws.Cells.Item(14, 5).Value = "NA"
$ws.Cells.Item(14, 6).Value = "NA"
$ws.Cells.Item(14, 7).Value = "W "
$ws.Cells.Item(15, 3).Value = "ARCH 252 and ARCH 242."
$ws.Cells.Item(15, 4).Value = "ARCH 207."
$ws.Cells.Item(15, 5).Value = "NA"
$ws.Cells.Item(15, 6).Value = "NA"
$ws.Cells.Item(15, 7).Value = "SP "
$ws.Cells.Item(16, 4).Value = "NA"
$ws.Cells.Item(16, 5).Value = "NA"
$ws.Cells.Item(16, 6).Value = "NA"
$ws.Cells.Item(16, 7).Value = "TBD"
$ws.Cells.Item(17, 4).Value = "NA"
$ws.Cells.Item(17, 5).Value = "NA"
$ws.Cells.Item(17, 6).Value = "NA"
$ws.Cells.Item(17, 7).Value = "TBD"
$ws.Cells.Item(18, 3).Value = "ARCH 341."
$ws.Cells.Item(18, 4).Value = "NA"
$ws.Cells.Item(18, 5).Value = "ARCH 352."
$ws.Cells.Item(18, 6).Value = "NA"
$ws.Cells.Item(18, 7).Value = "W "
$ws.Cells.Item(19, 4).Value = "NA"
$ws.Cells.Item(19, 5).Value = "NA"
$ws.Cells.Item(19, 6).Value = "NA"
$ws.Cells.Item(19, 7).Value = "SP"
$ws.Cells.Item(20, 3).Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of one course in GE Area B1 with a grade of C- or better; and one of the following GE Area C3 ARCH 217, 218, 219, or ART 112."
$ws.Cells.Item(20, 4).Value = "NA"
$ws.Cells.Item(20, 5).Value = "NA"
$ws.Cells.Item(20, 6).Value = "NA"
$ws.Cells.Item(20, 7).Value = "TBD"
$ws.Cells.Item(21, 4).Value = "NA"
$ws.Cells.Item(21, 5).Value = "NA"
$ws.Cells.Item(21, 6).Value = "NA"
$ws.Cells.Item(21, 7).Value = "W"
$ws.Cells.Item(22, 3).Value = "ARCH 207 and ARCH 253."
$ws.Cells.Item(22, 4).Value = "ARCH 351."
$ws.Cells.Item(22, 5).Value = "NA"
$ws.Cells.Item(22, 6).Value = "NA"
$ws.Cells.Item(22, 7).Value = "F "
$ws.Cells.Item(23, 3).Value = "ARCH 307."
$ws.Cells.Item(23, 4).Value = "NA"
$ws.Cells.Item(23, 5).Value = "ARCH 353."
$ws.Cells.Item(23, 6).Value = "NA"
$ws.Cells.Item(23, 7).Value = "SP "
$ws.Cells.Item(24, 3).Value = "ARCE 212, ARCH 253, ARCH 207 and PHYS 122 or PHYS 132, or consent of department head."
$ws.Cells.Item(24, 4).Value = "ARCH 341."
$ws.Cells.Item(24, 5).Value = "NA"
$ws.Cells.Item(24, 6).Value = "NA"
$ws.Cells.Item(24, 7).Value = "F "
$ws.Cells.Item(25, 3).Value = "ARCH 351, ARCH 341."
$ws.Cells.Item(25, 4).Value = "ARCH 307."
$ws.Cells.Item(25, 5).Value = "NA"
$ws.Cells.Item(25, 6).Value = "NA"
$ws.Cells.Item(25, 7).Value = "W "
$ws.Cells.Item(26, 3).Value = "ARCH 352, ARCH 307."
$ws.Cells.Item(26, 4).Value = "ARCH 342."
$ws.Cells.Item(26, 5).Value = "NA"
$ws.Cells.Item(26, 6).Value = "NA"
$ws.Cells.Item(26, 7).Value = "SP "
$ws.Cells.Item(27, 4).Value = "NA"
$ws.Cells.Item(27, 5).Value = "NA"
$ws.Cells.Item(27, 6).Value = "NA"
$ws.Cells.Item(27, 7).Value = "SP"
$ws.Cells.Item(28, 4).Value = "NA"
$ws.Cells.Item(28, 5).Value = "NA"
$ws.Cells.Item(28, 6).Value = "NA"
$ws.Cells.Item(28, 7).Value = "F,W,SP,SU"
$ws.Cells.Item(29, 4).Value = "NA"
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(29, 6).Value = "NA"
$ws.Cells.Item(29, 7).Value = "F, W, SP"
$ws.Cells.Item(30, 4).Value = "NA"
$ws.Cells.Item(30, 5).Value = "NA"
$ws.Cells.Item(30, 6).Value = "NA"
$ws.Cells.Item(30, 7).Value = "W"
$ws.Cells.Item(31, 4).Value = "NA"
$ws.Cells.Item(31, 5).Value = "NA"
$ws.Cells.Item(31, 6).Value = "NA"
$ws.Cells.Item(31, 7).Value = "TBD"
$ws.Cells.Item(32, 4).Value = "NA"
$ws.Cells.Item(32, 5).Value = "NA"
$ws.Cells.Item(32, 6).Value = "NA"
$ws.Cells.Item(32, 7).Value = "F"
$ws.Cells.Item(33, 4).Value = "NA"
$ws.Cells.Item(33, 5).Value = "NA"
$ws.Cells.Item(33, 6).Value = "NA"
$ws.Cells.Item(33, 7).Value = "W"
$ws.Cells.Item(34, 4).Value = "NA"
$ws.Cells.Item(34, 5).Value = "NA"
$ws.Cells.Item(34, 6).Value = "NA"
$ws.Cells.Item(34, 7).Value = "SP"
$ws.Cells.Item(35, 4).Value = "NA"
$ws.Cells.Item(35, 5).Value = "NA"
$ws.Cells.Item(35, 6).Value = "NA"
$ws.Cells.Item(35, 7).Value = "TBD"
$ws.Cells.Item(36, 4).Value = "NA"
$ws.Cells.Item(36, 5).Value = "NA"
$ws.Cells.Item(36, 6).Value = "NA"
$ws.Cells.Item(36, 7).Value = "TBD"
$ws.Cells.Item(37, 4).Value = "NA"
$ws.Cells.Item(37, 5).Value = "NA"
$ws.Cells.Item(37, 6).Value = "NA"
$ws.Cells.Item(37, 7).Value = "W"
$ws.Cells.Item(38, 4).Value = "NA"
$ws.Cells.Item(38, 5).Value = "NA"
$ws.Cells.Item(38, 6).Value = "NA"
$ws.Cells.Item(38, 7).Value = "TBD"
$ws.Cells.Item(39, 4).Value = "NA"
$ws.Cells.Item(39, 5).Value = "NA"
$ws.Cells.Item(39, 6).Value = "NA"
$ws.Cells.Item(39, 7).Value = "TBD"
$ws.Cells.Item(40, 4).Value = "NA"
$ws.Cells.Item(40, 5).Value = "NA"
$ws.Cells.Item(40, 6).Value = "NA"
$ws.Cells.Item(40, 7).Value = "TBD"
$ws.Cells.Item(41, 4).Value = "NA"
$ws.Cells.Item(41, 5).Value = "NA"
$ws.Cells.Item(41, 6).Value = "NA"
$ws.Cells.Item(41, 7).Value = "TBD"
$ws.Cells.Item(42, 4).Value = "NA"
$ws.Cells.Item(42, 5).Value = "NA"
$ws.Cells.Item(42, 6).Value = "NA"
$ws.Cells.Item(42, 7).Value = "TBD"
$ws.Cells.Item(43, 4).Value = "NA"
$ws.Cells.Item(43, 5).Value = "NA"
$ws.Cells.Item(43, 6).Value = "NA"
$ws.Cells.Item(43, 7).Value = "TBD"
$ws.Cells.Item(44, 4).Value = "NA"
$ws.Cells.Item(44, 5).Value = "NA"
$ws.Cells.Item(44, 6).Value = "NA"
$ws.Cells.Item(44, 7).Value = "TBD"
$ws.Cells.Item(45, 4).Value = "NA"
$ws.Cells.Item(45, 5).Value = "NA"
$ws.Cells.Item(45, 6).Value = "NA"
$ws.Cells.Item(45, 7).Value = "TBD"
$ws.Cells.Item(46, 4).Value = "NA"
$ws.Cells.Item(46, 5).Value = "NA"
$ws.Cells.Item(46, 6).Value = "NA"
$ws.Cells.Item(46, 7).Value = "TBD"
$ws.Cells.Item(47, 4).Value = "NA"
$ws.Cells.Item(47, 5).Value = "NA"
$ws.Cells.Item(47, 6).Value = "NA"
$ws.Cells.Item(47, 7).Value = "F,W,SP,SU"
$ws.Cells.Item(48, 4).Value = "NA"
$ws.Cells.Item(48, 5).Value = "NA"
$ws.Cells.Item(48, 6).Value = "NA"
$ws.Cells.Item(48, 7).Value = "F, W, SP"
$ws.Cells.Item(49, 4).Value = "NA"
$ws.Cells.Item(49, 5).Value = "NA"
$ws.Cells.Item(49, 6).Value = "NA"
$ws.Cells.Item(49, 7).Value = "TBD"
$ws.Cells.Item(50, 4).Value = "NA"
$ws.Cells.Item(50, 5).Value = "NA"
$ws.Cells.Item(50, 6).Value = "NA"
$ws.Cells.Item(50, 7).Value = "F,W,SP,SU"
$ws.Cells.Item(51, 3).Value = "ARCH 451, ARCH 452 and ARCH 453."
$ws.Cells.Item(51, 4).Value = "NA"
$ws.Cells.Item(51, 5).Value = "First quarter of ARCH 481."
$ws.Cells.Item(51, 6).Value = "NA"
$ws.Cells.Item(51, 7).Value = "F "
$ws.Cells.Item(52, 4).Value = "NA"
$ws.Cells.Item(52, 5).Value = "NA"
$ws.Cells.Item(52, 6).Value = "NA"
$ws.Cells.Item(52, 7).Value = "TBD"
$ws.Cells.Item(53, 4).Value = "NA"
$ws.Cells.Item(53, 5).Value = "NA"
$ws.Cells.Item(53, 6).Value = "NA"
$ws.Cells.Item(53, 7).Value = "F, W, SP"
$ws.Cells.Item(54, 4).Value = "NA"
$ws.Cells.Item(54, 5).Value = "NA"
$ws.Cells.Item(54, 6).Value = "NA"
$ws.Cells.Item(54, 7).Value = "F, W, SP"
$ws.Cells.Item(55, 4).Value = "NA"
$ws.Cells.Item(55, 5).Value = "NA"
$ws.Cells.Item(55, 6).Value = "NA"
$ws.Cells.Item(55, 7).Value = "TBD"
$ws.Cells.Item(56, 4).Value = "NA"
$ws.Cells.Item(56, 5).Value = "NA"
$ws.Cells.Item(56, 6).Value = "NA"
$ws.Cells.Item(56, 7).Value = "F, W, SP"
$ws.Cells.Item(57, 4).Value = "NA"
$ws.Cells.Item(57, 5).Value = "NA"
$ws.Cells.Item(57, 6).Value = "NA"
$ws.Cells.Item(57, 7).Value = "TBD"
